# Minor updates for docker handling and presentations
# The "Images" sheet had a stale "windmills-5643293_1280.jpg" / "#Klima #CO2 #Umweltschutz"
# row (row 3) removed; everything below shifts up. The Images sheet becomes the
# active/selected sheet with A2 selected, and Issues loses its "selected" flag.

$wb = $excel.ActiveWorkbook

$imagesSheet = $wb.Worksheets.Item("Images")
$imagesSheet.Activate()

$imagesSheet.Rows.Item(3).Delete()

$imagesSheet.Range("A2").Select()
